$wb = $excel.ActiveWorkbook

# --- TestData sheet: fill in K5/L5 automation id pair ---------------------
$wsData = $wb.Worksheets.Item("TestData")
$wsData.Range("L5").Value = "AUTOMATION1625826723452"

# --- TestSet sheet: insert a new "do not delete" execution-seed row -------
$wsSet = $wb.Worksheets.Item("TestSet")

# Insert a fresh row 2 (everything below shifts down one row).
$wsSet.Rows.Item(2).Insert()

# Seed its formatting from the row that is now directly below it (the old
# row 2, "LOGN_0001") so font/border match the rest of the table, then
# paint it yellow to flag it as the execution-seed row.
$srcRow = $wsSet.Range("A3:E3")
$srcRow.Copy()
$newRow = $wsSet.Range("A2:E2")
$newRow.PasteSpecial(-4122)
$newRow.Interior.Color = 65535

# Values for the new row (order matters so new shared strings land in the
# same sequence as the source workbook: L5 already added AUTOMATION..., so
# testExecution comes next, then the long "do not delete" label).
$wsSet.Range("E2").Value = "testExecution"
$wsSet.Range("A2").Value = "Donotdeleteusedforcreationoftestexecution"
$wsSet.Range("B2").Value = "YES"
$wsSet.Range("D2").Value = "LoginTest"
$wsSet.Range("C2").NumberFormat = "General"

# The two rows that used to be CRET_0001 / SEAR_0002 (now rows 4 & 5) lose
# their "YES" execution flag.
$wsSet.Range("B4").ClearContents()
$wsSet.Range("B5").ClearContents()

# Column A widens to fit the new long label.
$wsSet.Columns.Item(1).ColumnWidth = 40.71

# --- Selections / active sheet --------------------------------------------
$wsData.Range("A5").Select()
$wsSet.Activate()
$wsSet.Range("A3").Select()
